# Fixed Bento 80 Test scripts
# - CasesTab query (B2): append an ORDER BY / LIMIT clause
# - SamplesTab query (B3): append an ORDER BY / LIMIT clause
# - FilesTab query (B4): replace the trailing "order by" clause with a
#   capitalised ORDER BY / LIMIT clause
# - Row 3 grows taller (wrapped text got longer) -> bump its height
# - Scroll the view up one row (best effort; some hosts don't persist this)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SamplesTab query (cell B3): append order-by / limit ------------------
$samplesQuery = $ws.Range("B3").Value()
$samplesQuery = $samplesQuery + "  order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $samplesQuery

# --- CasesTab query (cell B2): append order-by / limit -------------------
$casesQuery = $ws.Range("B2").Value()
$casesQuery = $casesQuery + "  order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B2").Value = $casesQuery

# --- FilesTab query (cell B4): swap the lowercase order-by tail ----------
$filesQuery = $ws.Range("B4").Value()
$oldTail = "    order by f.file_name"
$newTail = " order By f.file_name ASC LIMIT 100"
$filesQuery = $filesQuery.Substring(0, $filesQuery.Length - $oldTail.Length) + $newTail
$ws.Range("B4").Value = $filesQuery

# --- Row 3 is now taller because its wrapped cells hold more text --------
$ws.Rows.Item(3).RowHeight = 360

# --- Scroll the window so row 3 is the top visible row -------------------
try {
    $ws.Activate()
    $excel.ActiveWindow.ScrollRow = 3
} catch {
}
